# Applies the changes described by the diff:
#  - "Reguły" sheet: reorder the object lists inside 4 of the 6 rule texts
#  - "Walidacja krzyżowa" sheet: reorder the metric rows (labels + values
#    move together, i.e. values stay attached to their original label)

$wb = $excel.ActiveWorkbook

# --- Sheet "Reguły" (rules) : column B holds the rule text ---------------
$wsRules = $wb.Worksheets.Item("Reguły")

$wsRules.Range("B2").Value = "(attempts >=  3.0) => (class <= 1) ['a23', 'a7', 'a13', 'a1', 'a3']"
$wsRules.Range("B3").Value = "(pregnancy <=  0.0) & (age >=  32.0) & (frozen_embryos <=  4.0) & (cleavage_stage >=  5.0) => (class <= 1) ['a15', 'a21', 'a22', 'a2', 'a3']"
$wsRules.Range("B5").Value = "(age >=  42.0) => (class <= 1) ['a3', 'a14']"
$wsRules.Range("B6").Value = "(age <=  31.0) & (attempts <=  1.0) => (class >= 2) ['a12', 'a25', 'a24', 'a11', 'a9']"
$wsRules.Range("B7").Value = "(frozen_embryos >=  8.0) & (sperm <=  1.0) => (class >= 2) ['a6', 'a16']"

# --- Sheet "Walidacja krzyżowa" (cross validation) : reorder rows --------
$wsCv = $wb.Worksheets.Item("Walidacja krzyżowa")

$wsCv.Range("A1").Value = "correct"
$wsCv.Range("B1").Value = 0.7142857142857143

$wsCv.Range("A2").Value = "not_classified"
$wsCv.Range("B2").Value = 0.44

$wsCv.Range("A3").Value = "f1_score"
$wsCv.Range("B3").Value = 0.4952380952380952

$wsCv.Range("A4").Value = "accuracy"
$wsCv.Range("B4").Value = 0.4
